$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (style uses date format).
# Update rows 2-119 from 45182 (2023-09-13) to 45184 (2023-09-15).
$ws.Range("C2:C119").Value = 45184
